$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "2026-02-20T08:03:27.271167+00:00"
$ws.Range("H16").Value = 3
$ws.Range("L16").Value = "[101305, 101325, 101324]"
